# Apply the "operator change" edit described by the commit:
#   - Car GL757TH moves from operator GIAMPIERO.RICUCCI to DI LANZO ALESSIO
#     on 2026-02-12, recorded on the "Stato Attuale" sheet.
#   - The "Storico Passaggi" (change history) sheet is updated to contain
#     only this single new change-of-operator entry.

$wb = $excel.ActiveWorkbook

# ---- Sheet "Stato Attuale": update the current operator for GL757TH ----
$ws1 = $wb.Worksheets.Item("Stato Attuale")

$ws1.Range("B82").Value = "DI LANZO ALESSIO"

# Write the date as plain text (not an auto-converted date serial), same as
# every other date cell in this column: enter with a leading apostrophe so
# it is stored as text, then drop the resulting "text" number format so the
# cell keeps the sheet's default (unstyled) formatting.
$ws1.Range("C82").Value = "'2026-02-12"
$ws1.Range("C82").ClearFormats()

# ---- Sheet "Storico Passaggi": replace history with the new entry ----
$ws2 = $wb.Worksheets.Item("Storico Passaggi")

# Drop the old history rows (rows 3 and 4); row 2 will be overwritten below.
$ws2.Range("A3:D4").ClearContents()

$ws2.Range("A2").Value = "GL757TH"
$ws2.Range("B2").Value = "GIAMPIERO.RICUCCI"
$ws2.Range("C2").Value = "DI LANZO ALESSIO"
$ws2.Range("D2").Value = "'2026-02-12"
$ws2.Range("D2").ClearFormats()
